# =====================================================================
# Edit script: "Unraveling the Cosmic Tapestry" (black holes) ->
#              "The Magical World of Chemistry"
# =====================================================================
$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

# --- Title -------------------------------------------------------------
Replace-Text "Unraveling the Cosmic Tapestry: The Symphony of Black Holes" "The Magical World of Chemistry"

# --- Author name: "Nicholas Corrigan" -> "Dr" + "." + " Bethany Isabelle" (3 runs) ---
$rng = $d.Content
$rng.Find.Execute("Nicholas Corrigan") | Out-Null
$rng.Text = "Dr"
$afterDr = $rng.Duplicate
$afterDr.Collapse(0)
$afterDr.InsertAfter(".")
$afterDr.Bold = 1
$afterDr.Bold = 0

$afterDot = $afterDr.Duplicate
$afterDot.Collapse(0)
$afterDot.InsertAfter(" Bethany Isabelle")
$afterDot.Bold = 1
$afterDot.Bold = 0

# --- Email -------------------------------------------------------------
Replace-Text "ncorrigan@stargazeobservatory" "bkisabelle@libertyview"
Replace-Text "org" "com"

# --- Body paragraph 1 ----------------------------------------------------
Replace-Text "Emerging from the depths of cosmic darkness, black holes stand as enigmatic entities, captivating the imaginations of scientists and philosophers alike" "In the captivating world of chemistry, every moment is a thrilling revelation of the universe's most fundamental secrets"

Replace-Text " These celestial behemoths, born from the gravitational collapse of massive stars, defy conventional logic and test the limits of human comprehension" " It is a realm where tiny atoms come together in a grand cosmic dance, forming the very structure of everything around us, from the air we breathe to the stars twinkling in the night sky"

Replace-Text " In this cosmic symphony, black holes conduct a mesmerizing dance, orchestrating the fate of matter and energy within their gravitational grasp" " Chemistry unveils the intricacies of nature, revealing the mystical ballet of electrons as they waltz around the nucleus, orchestrated by the universal rules of attraction and repulsion"

# Delete "." + " Delving into the realm of astrophysics..." runs entirely
$delRng = $d.Content
$delRng.Find.Execute(". Delving into the realm of astrophysics, we embark on a quest to unravel the mysteries surrounding these enigmatic cosmic wonders") | Out-Null
$delRng.Text = ""

# --- "Unveiling the Secrets of Gravity's Embrace:" heading line --------
Replace-Text "Unveiling the Secrets of Gravity's Embrace:" "This mesmerizing show is performed everywhere, from the intricate workings of our own bodies to the majestic processes of the cosmos"

# Insert two new runs after it: "." and " In this grand tapestry..."
$afterShow = $d.Content
$afterShow.Find.Execute("This mesmerizing show is performed everywhere, from the intricate workings of our own bodies to the majestic processes of the cosmos") | Out-Null
$afterShow.Collapse(0)
$afterShow.InsertAfter(".")
$afterShow.Bold = 1
$afterShow.Bold = 0

$afterShowDot = $afterShow.Duplicate
$afterShowDot.Collapse(0)
$afterShowDot.InsertAfter(" In this grand tapestry of chemistry, we find the answers to questions that have perplexed humanity for ages: Why do leaves change color in the fall? How do stars generate their magnificent energy?  What makes a certain food taste the way it does?")
$afterShowDot.Bold = 1
$afterShowDot.Bold = 0

# --- "Gravity, the invisible thread..." paragraph heading line ---------
Replace-Text "Gravity, the invisible thread that weaves the fabric of the universe, exerts its dominion over black holes like a celestial maestro" "As we embark on this extraordinary exploration of chemistry, we will embark on an odyssey through the ages, uncovering the pioneering contributions made by esteemed scientists who dedicated their lives to unraveling the secrets of matter"

# --- Big deletion: from " Within their event horizons..." through
#     " Additionally, the principles governing black holes...reach", replaced
#     by a single new sentence.
$bigRng = $d.Content
$bigStart = $d.Content
$bigStart.Find.Execute(" Within their event horizons") | Out-Null
$bigEnd = $d.Content
$bigEnd.Find.Execute("Additionally, the principles governing black holes could inspire innovations in gravitational wave detection and space exploration, propelling humankind towards a future where the mysteries of the cosmos are brought within our reach") | Out-Null
$combined = $d.Range($bigStart.Start, $bigEnd.End)
$combined.Text = " Their discoveries have transformed our understanding of the universe, leading to advancements in medicine, energy, and technology that shape our modern-day existence"

# --- Summary section -----------------------------------------------------
Replace-Text "Black holes, cosmic entities of immense gravitational pull and shrouded in mystery, beckon us to explore the profound interplay of gravity, singularity, and the information paradox" "Chemistry, the dazzling dance of elements, offers an entrancing glimpse into the universe's fundamental mysteries"

Replace-Text " By unraveling the secrets of these cosmic behemoths, we not only deepen our understanding of the universe but also pave the way for groundbreaking technological advancements" " It unravels the secrets of nature, elucidates the transformation of substances, and chronicles the inventive minds that propelled scientific progress"

Replace-Text " The symphony of black holes, with its captivating blend of enigma and potential, will continue to inspire awe and fuel our scientific quest for knowledge for generations to come" " From the microscopic interactions of molecules to the mesmerizing spectacles of celestial bodies, chemistry offers a fascinating inquiry into the very fabric of our existence"

# --- Add a trailing empty paragraph at the very end of the document ------
$endRng = $d.Content
$endRng.Collapse(0)
$endRng.InsertParagraphAfter()
